$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for the summary block (J1:O1)
$ws.Range("J1").Value = "销售地区"
$ws.Range("K1").Value = "销售一分部"
$ws.Range("L1").Value = "销售三分部"
$ws.Range("M1").Value = "销售二分部"
$ws.Range("N1").Value = "销售四分部"
$ws.Range("O1").Value = "总计"

# Row 2 - 华东
$ws.Range("J2").Value = "华东"
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 2059200
$ws.Range("M2").Value = 4183800
$ws.Range("N2").Value = 3513200
$ws.Range("O2").Value = 9756200

# Row 3 - 华中
$ws.Range("J3").Value = "华中"
$ws.Range("K3").Value = 3826000
$ws.Range("L3").Value = 1806200
$ws.Range("M3").Value = 6324000
$ws.Range("N3").Value = 1284000
$ws.Range("O3").Value = 13240200

# Row 4 - 华北
$ws.Range("J4").Value = "华北"
$ws.Range("K4").Value = 3676400
$ws.Range("L4").Value = 1694000
$ws.Range("M4").Value = 1245200
$ws.Range("N4").Value = 2552800
$ws.Range("O4").Value = 9168400

# Row 5 - 华南
$ws.Range("J5").Value = "华南"
$ws.Range("K5").Value = 3025200
$ws.Range("L5").Value = 1634600
$ws.Range("M5").Value = 588000
$ws.Range("N5").Value = 3369400
$ws.Range("O5").Value = 8617200

# Row 6 - 总计
$ws.Range("J6").Value = "总计"
$ws.Range("K6").Value = 10527600
$ws.Range("L6").Value = 7194000
$ws.Range("M6").Value = 12341000
$ws.Range("N6").Value = 10719400
$ws.Range("O6").Value = 40782000

# Apply currency number format (with 2 decimals) to the numeric block
$ws.Range("L2:O6").NumberFormat = "¥#,##0.00;¥-#,##0.00"
$ws.Range("K3:K6").NumberFormat = "¥#,##0.00;¥-#,##0.00"
$ws.Range("J1:O6").HorizontalAlignment = -4108

$ws.Columns("A:H").AutoFit()
$ws.Columns("J:O").AutoFit()

$ws.Range("K14").Select()
